# edit.ps1 - Updates CFR_Results.xlsx "Results" sheet for R2018b/CFL v4.1 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): software version stamps + run date/time ---
# H1 is set before C1 so the shared-string table allocates indices in the
# same order as the reference edit (CFL version string before MATLAB string).
$ws.Range("H1").Value = "CFL v4.1"
$ws.Range("C1").Value = "9.5.0.944444 (R2018b)"
$ws.Range("F1").Value = 43354.271284722221

# --- Restore the frozen-pane scroll position to the top of the data table ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 3

# --- Refreshed benchmark results (rows 4-55) ---
# Row 4
$ws.Cells.Item(4,5).Value = 0.98899999999999999
$ws.Cells.Item(4,7).Value = 0.16
$ws.Cells.Item(4,9).Value = 0.13800000000000001
$ws.Cells.Item(4,11).Value = 0.155
# Row 5
$ws.Cells.Item(5,5).Value = 0.63600000000000001
$ws.Cells.Item(5,6).Value = 935
$ws.Cells.Item(5,7).Value = 0.107
$ws.Cells.Item(5,9).Value = 0.124
$ws.Cells.Item(5,11).Value = 0.107
# Row 6
$ws.Cells.Item(6,5).Value = 0.23499999999999999
$ws.Cells.Item(6,7).Value = 0.077
$ws.Cells.Item(6,9).Value = 0.071
$ws.Cells.Item(6,11).Value = 0.088
# Row 7
$ws.Cells.Item(7,5).Value = 1.871
$ws.Cells.Item(7,7).Value = 2.57
$ws.Cells.Item(7,9).Value = 2.5569999999999999
$ws.Cells.Item(7,11).Value = 3.52
# Row 8
$ws.Cells.Item(8,4).Value = 2071
$ws.Cells.Item(8,5).Value = 0.29199999999999998
$ws.Cells.Item(8,6).Value = 2660
$ws.Cells.Item(8,7).Value = 0.33800000000000002
$ws.Cells.Item(8,8).Value = 2276
$ws.Cells.Item(8,9).Value = 0.23799999999999999
$ws.Cells.Item(8,10).Value = 2954
$ws.Cells.Item(8,11).Value = 0.33200000000000002
# Row 9
$ws.Cells.Item(9,5).Value = 0.22600000000000001
$ws.Cells.Item(9,7).Value = 0.12
$ws.Cells.Item(9,9).Value = 0.13300000000000001
$ws.Cells.Item(9,11).Value = 0.14199999999999999
# Row 10
$ws.Cells.Item(10,5).Value = 0.159
$ws.Cells.Item(10,7).Value = 0.052
$ws.Cells.Item(10,9).Value = 0.051
$ws.Cells.Item(10,11).Value = 0.051
# Row 11
$ws.Cells.Item(11,5).Value = 2.633
$ws.Cells.Item(11,7).Value = 5.1280000000000001
$ws.Cells.Item(11,9).Value = 3.343
$ws.Cells.Item(11,11).Value = 6.5
# Row 12
$ws.Cells.Item(12,5).Value = 0.10199999999999999
$ws.Cells.Item(12,7).Value = 0.115
$ws.Cells.Item(12,9).Value = 0.063
$ws.Cells.Item(12,11).Value = 0.075
# Row 13
$ws.Cells.Item(13,5).Value = 0.255
$ws.Cells.Item(13,7).Value = 0.19
$ws.Cells.Item(13,9).Value = 0.16900000000000001
$ws.Cells.Item(13,11).Value = 0.26300000000000001
# Row 14
$ws.Cells.Item(14,5).Value = 0.308
$ws.Cells.Item(14,7).Value = 0.41899999999999998
$ws.Cells.Item(14,9).Value = 0.35699999999999998
$ws.Cells.Item(14,11).Value = 0.377
# Row 15
$ws.Cells.Item(15,5).Value = 0.184
$ws.Cells.Item(15,7).Value = 0.085
$ws.Cells.Item(15,9).Value = 0.082
$ws.Cells.Item(15,10).Value = 989
$ws.Cells.Item(15,11).Value = 0.105
# Row 16
$ws.Cells.Item(16,5).Value = 0.185
$ws.Cells.Item(16,7).Value = 0.17199999999999999
$ws.Cells.Item(16,9).Value = 0.111
$ws.Cells.Item(16,11).Value = 0.151
# Row 17
$ws.Cells.Item(17,5).Value = 1.2929999999999999
$ws.Cells.Item(17,7).Value = 1.135
$ws.Cells.Item(17,9).Value = 1.079
$ws.Cells.Item(17,11).Value = 1.2789999999999999
# Row 18
$ws.Cells.Item(18,5).Value = 0.13900000000000001
$ws.Cells.Item(18,7).Value = 0.066
$ws.Cells.Item(18,9).Value = 0.071
$ws.Cells.Item(18,11).Value = 0.069
# Row 19
$ws.Cells.Item(19,4).Value = 3443
$ws.Cells.Item(19,5).Value = 0.747
$ws.Cells.Item(19,6).Value = 3675
$ws.Cells.Item(19,7).Value = 0.63900000000000001
$ws.Cells.Item(19,8).Value = 3496
$ws.Cells.Item(19,9).Value = 0.52800000000000002
$ws.Cells.Item(19,10).Value = 3803
$ws.Cells.Item(19,11).Value = 0.68799999999999994
# Row 20
$ws.Cells.Item(20,4).Value = 6874
$ws.Cells.Item(20,5).Value = 0.83399999999999996
$ws.Cells.Item(20,6).Value = 7094
$ws.Cells.Item(20,7).Value = 0.88300000000000001
$ws.Cells.Item(20,8).Value = 98871
$ws.Cells.Item(20,9).Value = 9.9429999999999996
$ws.Cells.Item(20,10).Value = 99867
$ws.Cells.Item(20,11).Value = 11.047000000000001
# Row 21
$ws.Cells.Item(21,4).Value = 3697
$ws.Cells.Item(21,5).Value = 0.63400000000000001
$ws.Cells.Item(21,6).Value = 4605
$ws.Cells.Item(21,7).Value = 0.53900000000000003
$ws.Cells.Item(21,8).Value = 3969
$ws.Cells.Item(21,9).Value = 0.375
$ws.Cells.Item(21,10).Value = 3933
$ws.Cells.Item(21,11).Value = 0.38400000000000001
# Row 22
$ws.Cells.Item(22,4).Value = 2229
$ws.Cells.Item(22,5).Value = 4.8159999999999998
$ws.Cells.Item(22,6).Value = 2229
$ws.Cells.Item(22,7).Value = 4.3360000000000003
$ws.Cells.Item(22,8).Value = 2229
$ws.Cells.Item(22,9).Value = 4.242
$ws.Cells.Item(22,10).Value = 2229
$ws.Cells.Item(22,11).Value = 4.2729999999999997
# Row 23
$ws.Cells.Item(23,5).Value = 1.262
$ws.Cells.Item(23,7).Value = 1.256
$ws.Cells.Item(23,9).Value = 1.2869999999999999
$ws.Cells.Item(23,11).Value = 1.2529999999999999
# Row 24
$ws.Cells.Item(24,5).Value = 0.14299999999999999
$ws.Cells.Item(24,7).Value = 0.053
$ws.Cells.Item(24,9).Value = 0.056
$ws.Cells.Item(24,11).Value = 0.059
# Row 25
$ws.Cells.Item(25,5).Value = 3.8479999999999999
$ws.Cells.Item(25,7).Value = 4.2699999999999996
$ws.Cells.Item(25,9).Value = 3.871
$ws.Cells.Item(25,11).Value = 4.5529999999999999
# Row 26
$ws.Cells.Item(26,4).Value = 14171
$ws.Cells.Item(26,5).Value = 6.2629999999999999
$ws.Cells.Item(26,6).Value = 14108
$ws.Cells.Item(26,7).Value = 6.0140000000000002
$ws.Cells.Item(26,8).Value = 20521
$ws.Cells.Item(26,9).Value = 8.7970000000000006
$ws.Cells.Item(26,10).Value = 19890
$ws.Cells.Item(26,11).Value = 8.6489999999999991
# Row 27
$ws.Cells.Item(27,5).Value = 2.2690000000000001
$ws.Cells.Item(27,7).Value = 3.4060000000000001
$ws.Cells.Item(27,9).Value = 2.1539999999999999
$ws.Cells.Item(27,11).Value = 3.4159999999999999
# Row 28
$ws.Cells.Item(28,5).Value = 0.249
$ws.Cells.Item(28,7).Value = 0.16800000000000001
$ws.Cells.Item(28,9).Value = 0.13100000000000001
$ws.Cells.Item(28,11).Value = 0.156
# Row 29
$ws.Cells.Item(29,5).Value = 0.21299999999999999
$ws.Cells.Item(29,9).Value = 0.221
$ws.Cells.Item(29,11).Value = 0.216
# Row 30
$ws.Cells.Item(30,5).Value = 12.255000000000001
$ws.Cells.Item(30,7).Value = 24.074000000000002
$ws.Cells.Item(30,9).Value = 18.384
$ws.Cells.Item(30,11).Value = 31.193000000000001
# Row 31
$ws.Cells.Item(31,5).Value = 0.41599999999999998
$ws.Cells.Item(31,7).Value = 0.33600000000000002
$ws.Cells.Item(31,9).Value = 0.25
$ws.Cells.Item(31,11).Value = 0.27800000000000002
# Row 32
$ws.Cells.Item(32,5).Value = 1.9339999999999999
$ws.Cells.Item(32,7).Value = 8.1999999999999993
$ws.Cells.Item(32,9).Value = 1.3360000000000001
$ws.Cells.Item(32,11).Value = 9.2289999999999992
# Row 33
$ws.Cells.Item(33,5).Value = 0.58599999999999997
$ws.Cells.Item(33,7).Value = 0.85299999999999998
$ws.Cells.Item(33,9).Value = 0.49399999999999999
$ws.Cells.Item(33,11).Value = 0.88
# Row 34
$ws.Cells.Item(34,5).Value = 0.221
$ws.Cells.Item(34,7).Value = 0.40799999999999997
$ws.Cells.Item(34,9).Value = 0.127
$ws.Cells.Item(34,11).Value = 0.41299999999999998
# Row 35
$ws.Cells.Item(35,5).Value = 0.318
$ws.Cells.Item(35,7).Value = 0.872
$ws.Cells.Item(35,9).Value = 0.23799999999999999
$ws.Cells.Item(35,11).Value = 0.876
# Row 36
$ws.Cells.Item(36,5).Value = 0.93899999999999995
$ws.Cells.Item(36,7).Value = 2.0049999999999999
$ws.Cells.Item(36,9).Value = 0.749
$ws.Cells.Item(36,11).Value = 2.2719999999999998
# Row 37
$ws.Cells.Item(37,5).Value = 0.85199999999999998
$ws.Cells.Item(37,7).Value = 2.2879999999999998
$ws.Cells.Item(37,9).Value = 0.64900000000000002
$ws.Cells.Item(37,11).Value = 2.444
# Row 38
$ws.Cells.Item(38,5).Value = 1.071
$ws.Cells.Item(38,7).Value = 1.42
$ws.Cells.Item(38,9).Value = 1.044
$ws.Cells.Item(38,11).Value = 1.4219999999999999
# Row 39
$ws.Cells.Item(39,5).Value = 0.64100000000000001
$ws.Cells.Item(39,7).Value = 0.68100000000000005
$ws.Cells.Item(39,9).Value = 0.54
$ws.Cells.Item(39,11).Value = 0.67700000000000005
# Row 40
$ws.Cells.Item(40,5).Value = 1.677
$ws.Cells.Item(40,7).Value = 2.0859999999999999
$ws.Cells.Item(40,9).Value = 1.6439999999999999
$ws.Cells.Item(40,11).Value = 2.0009999999999999
# Row 41
$ws.Cells.Item(41,5).Value = 1.373
$ws.Cells.Item(41,7).Value = 1.407
$ws.Cells.Item(41,9).Value = 1.3280000000000001
$ws.Cells.Item(41,11).Value = 1.5640000000000001
# Row 42
$ws.Cells.Item(42,5).Value = 0.71799999999999997
$ws.Cells.Item(42,7).Value = 0.86599999999999999
$ws.Cells.Item(42,9).Value = 0.73399999999999999
$ws.Cells.Item(42,11).Value = 0.93500000000000005
# Row 43
$ws.Cells.Item(43,5).Value = 0.53900000000000003
$ws.Cells.Item(43,7).Value = 0.45400000000000001
$ws.Cells.Item(43,9).Value = 0.27
$ws.Cells.Item(43,11).Value = 0.42599999999999999
# Row 44
$ws.Cells.Item(44,5).Value = 4.0890000000000004
$ws.Cells.Item(44,7).Value = 5.3159999999999998
$ws.Cells.Item(44,9).Value = 3.81
$ws.Cells.Item(44,11).Value = 5.3440000000000003
# Row 45
$ws.Cells.Item(45,4).Value = 4242
$ws.Cells.Item(45,5).Value = 2.9809999999999999
$ws.Cells.Item(45,6).Value = 4242
$ws.Cells.Item(45,7).Value = 2.5579999999999998
$ws.Cells.Item(45,8).Value = 4242
$ws.Cells.Item(45,9).Value = 2.9809999999999999
$ws.Cells.Item(45,10).Value = 4242
$ws.Cells.Item(45,11).Value = 2.5710000000000002
# Row 46
$ws.Cells.Item(46,3).Value = 50
$ws.Cells.Item(46,4).Value = 5515
$ws.Cells.Item(46,5).Value = 3.7170000000000001
$ws.Cells.Item(46,6).Value = 9462
$ws.Cells.Item(46,7).Value = 6.9489999999999998
$ws.Cells.Item(46,8).Value = 5558
$ws.Cells.Item(46,9).Value = 3.22
$ws.Cells.Item(46,10).Value = 9667
$ws.Cells.Item(46,11).Value = 7.3250000000000002
# Row 47
$ws.Cells.Item(47,4).Value = 9375
$ws.Cells.Item(47,5).Value = 5.2859999999999996
# Row 48
$ws.Cells.Item(48,4).Value = 8828
$ws.Cells.Item(48,5).Value = 5.0209999999999999
# Row 49
$ws.Cells.Item(49,4).Value = 15140
$ws.Cells.Item(49,5).Value = 11.003
# Row 50
$ws.Cells.Item(50,3).Value = 39.020000000000003
$ws.Cells.Item(50,4).Value = 11499
$ws.Cells.Item(50,5).Value = 7.7409999999999997
# Row 51
$ws.Cells.Item(51,3).Value = 34.68
$ws.Cells.Item(51,4).Value = 9443
$ws.Cells.Item(51,5).Value = 5.1580000000000004
# Row 52
$ws.Cells.Item(52,4).Value = 11825
$ws.Cells.Item(52,5).Value = 11.680999999999999
# Row 53
$ws.Cells.Item(53,3).Value = 41.28
$ws.Cells.Item(53,4).Value = 11517
$ws.Cells.Item(53,5).Value = 10.426
# Row 54
$ws.Cells.Item(54,3).Value = 45.27
$ws.Cells.Item(54,4).Value = 12245
$ws.Cells.Item(54,5).Value = 14.163
# Row 55
$ws.Cells.Item(55,4).Value = 8945
$ws.Cells.Item(55,5).Value = 8.5510000000000002
